$d = $word.ActiveDocument

# Locate the three paragraphs that must be removed:
#   1) the empty paragraph right after "LOQ4205: Sistemas Produtivos II (Requisito fraco)"
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "Contact: luizeleno@usp.br..." footer paragraph
# Anchor on the "LOQ4205" paragraph, then walk forward using .Next so we
# don't depend on matching special characters (e.g. the copyright sign).

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOQ4205: Sistemas Produtivos II*") {
        $anchor = $p
        break
    }
}

$p1 = $anchor.Next()
$p2 = $p1.Next()
$p3 = $p2.Next()

$p3.Range.Delete()
$p2.Range.Delete()
$p1.Range.Delete()
